{"js": "// Office.js (Word JavaScript API) script.\n// Applies the text replacements described by the diff: the header date\n// and each \"NNN\u00f7N=\" division prompt in the table, in document order.\n\nconst replacements = [\n  [\"2024-08-27 Tuesday\", \"2024-08-28 Wednesday\"],\n  [\"515\u00f79=\", \"856\u00f75=\"],\n  [\"910\u00f74=\", \"156\u00f75=\"],\n  [\"558\u00f73=\", \"154\u00f73=\"],\n  [\"634\u00f72=\", \"775\u00f78=\"],\n  [\"973\u00f73=\", \"954\u00f73=\"],\n  [\"697\u00f73=\", \"155\u00f76=\"],\n  [\"918\u00f75=\", \"484\u00f73=\"],\n  [\"312\u00f74=\", \"901\u00f74=\"],\n  [\"875\u00f73=\", \"817\u00f75=\"],\n  [\"174\u00f72=\", \"110\u00f77=\"],\n  [\"463\u00f76=\", \"229\u00f75=\"],\n  [\"350\u00f73=\", \"920\u00f75=\"],\n  [\"578\u00f74=\", \"649\u00f77=\"],\n  [\"421\u00f72=\", \"384\u00f78=\"],\n  [\"732\u00f75=\", \"235\u00f78=\"],\n  [\"678\u00f76=\", \"971\u00f79=\"],\n  [\"218\u00f72=\", \"557\u00f74=\"],\n  [\"639\u00f72=\", \"951\u00f77=\"],\n  [\"306\u00f76=\", \"408\u00f76=\"],\n  [\"827\u00f76=\", \"688\u00f79=\"],\n  [\"214\u00f79=\", \"428\u00f79=\"],\n  [\"878\u00f72=\", \"943\u00f74=\"],\n  [\"401\u00f73=\", \"861\u00f78=\"],\n  [\"823\u00f76=\", \"543\u00f79=\"],\n  [\"430\u00f76=\", \"274\u00f74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the text replacements described by the diff: the header date\n# and each \"NNN\u00f7N=\" division prompt in the table, in document order.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-08-27 Tuesday\", \"2024-08-28 Wednesday\"),\n    @(\"515\u00f79=\", \"856\u00f75=\"),\n    @(\"910\u00f74=\", \"156\u00f75=\"),\n    @(\"558\u00f73=\", \"154\u00f73=\"),\n    @(\"634\u00f72=\", \"775\u00f78=\"),\n    @(\"973\u00f73=\", \"954\u00f73=\"),\n    @(\"697\u00f73=\", \"155\u00f76=\"),\n    @(\"918\u00f75=\", \"484\u00f73=\"),\n    @(\"312\u00f74=\", \"901\u00f74=\"),\n    @(\"875\u00f73=\", \"817\u00f75=\"),\n    @(\"174\u00f72=\", \"110\u00f77=\"),\n    @(\"463\u00f76=\", \"229\u00f75=\"),\n    @(\"350\u00f73=\", \"920\u00f75=\"),\n    @(\"578\u00f74=\", \"649\u00f77=\"),\n    @(\"421\u00f72=\", \"384\u00f78=\"),\n    @(\"732\u00f75=\", \"235\u00f78=\"),\n    @(\"678\u00f76=\", \"971\u00f79=\"),\n    @(\"218\u00f72=\", \"557\u00f74=\"),\n    @(\"639\u00f72=\", \"951\u00f77=\"),\n    @(\"306\u00f76=\", \"408\u00f76=\"),\n    @(\"827\u00f76=\", \"688\u00f79=\"),\n    @(\"214\u00f79=\", \"428\u00f79=\"),\n    @(\"878\u00f72=\", \"943\u00f74=\"),\n    @(\"401\u00f73=\", \"861\u00f78=\"),\n    @(\"823\u00f76=\", \"543\u00f79=\"),\n    @(\"430\u00f76=\", \"274\u00f74=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
